{"js": "// Change the quoted filename \"Mixed_Ref.csv\" to \"Mixed_Ref.tsv\", splitting the\n// text into the same run/proofErr structure Word produces when it spell-checks\n// a freshly (re)typed word:\n//   \"  |  Mixed_Ref  |  .tsv  |  \"\n// with a spellStart/spellEnd proofErr pair bracketing \"Mixed_Ref\" + \".tsv\".\n\nconst OPEN_QUOTE = \"\\u201C\"; // \u201c\nconst CLOSE_QUOTE = \"\\u201D\"; // \u201d\nconst OLD_TEXT = OPEN_QUOTE + \"Mixed_Ref.csv\" + CLOSE_QUOTE;\n\nconst body = context.document.body;\nconst results = body.search(OLD_TEXT, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n    const match = results.items[0];\n\n    // 1) Remove the old text first (as its own sync) so the freshly inserted\n    //    runs below don't get merged with whatever formatting-compatible run\n    //    used to sit next to the deleted text.\n    match.insertText(\"\", Word.InsertLocation.replace);\n    await context.sync();\n\n    // 2) Insert the replacement as raw OOXML so we control the exact run\n    //    boundaries and the w:proofErr markers, matching what Word's editor\n    //    leaves behind after an in-place spelling-triggered re-split.\n    const ooxml =\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' +\n        '<w:r><w:t>' + OPEN_QUOTE + '</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>Mixed_Ref</w:t></w:r>' +\n        '<w:r><w:t>.tsv</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t>' + CLOSE_QUOTE + '</w:t></w:r>' +\n        '</w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>';\n\n    match.insertOoxml(ooxml, Word.InsertLocation.before);\n    await context.sync();\n}\n", "ps1": "# Change the quoted filename \"Mixed_Ref.csv\" to \"Mixed_Ref.tsv\", splitting the\n# text into the same run/proofErr structure Word produces when it spell-checks\n# a freshly (re)typed word:\n#   \"  |  Mixed_Ref  |  .tsv  |  \"\n# with a spellStart/spellEnd proofErr pair bracketing \"Mixed_Ref\" + \".tsv\".\n\n$d = $word.ActiveDocument\n\n$openQuote = [char]0x201C\n$closeQuote = [char]0x201D\n$oldText = $openQuote + \"Mixed_Ref.csv\" + $closeQuote\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = $oldText\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    $insertStart = $rng.Start\n\n    # 1) Delete the old text first so the new runs inserted below don't get\n    #    merged into whatever run used to be adjacent to it.\n    $rng.Delete()\n\n    # 2) Insert the replacement as raw OOXML at the now-collapsed point, so we\n    #    control the exact run boundaries and the w:proofErr markers.\n    $insertPoint = $d.Range($insertStart, $insertStart)\n    $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:r><w:t>' + $openQuote + '</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/>' +\n        '<w:r><w:t>Mixed_Ref</w:t></w:r>' +\n        '<w:r><w:t>.tsv</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t>' + $closeQuote + '</w:t></w:r>' +\n        '</w:p>'\n    $insertPoint.InsertXML($xml)\n}\n"}
